$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: copy formatting (style index 1 = bold/border/center) onto all
# destination cells that need it, while source cells still carry that style.
$ws.Range("B1").Copy()
$ws.Range("A1").PasteSpecial(-4122)

$ws.Range("A2").Copy()
$ws.Range("E1").PasteSpecial(-4122)

# C1 already has style 1 and stays in place; D1 also already has style 1 and
# stays styled (value changes only). B1 keeps style 1 too.
# A2..A5 already carry style 1 and remain in place.

$excel.CutCopyMode = 0

# --- Step 2: set the header row values (row 1)
$ws.Range("A1").Value = "Criterium"
$ws.Range("B1").Value = "mean"
$ws.Range("C1").Value = "max"
$ws.Range("D1").Value = "min"
$ws.Range("E1").Value = "meteomarge"

# --- Step 3: set row labels (column A) and numeric data (columns B-E)
$ws.Range("A2").Value = "Won 58 dB(A) Lden"
$ws.Range("B2").Value = 10200
$ws.Range("C2").Value = 10800
$ws.Range("D2").Value = 9700
$ws.Range("E2").Value = 11300

$ws.Range("A3").Value = "EGH 48 dB(A) Lden"
$ws.Range("B3").Value = 98800
$ws.Range("C3").Value = 103000
$ws.Range("D3").Value = 94800
$ws.Range("E3").Value = 112000

$ws.Range("A4").Value = "Won 48 dB(A) Lnight"
$ws.Range("B4").Value = 7300
$ws.Range("C4").Value = 8400
$ws.Range("D4").Value = 6000
$ws.Range("E4").Value = 9000

$ws.Range("A5").Value = "SV 40 dB(A) Lnight"
$ws.Range("B5").Value = 20600
$ws.Range("C5").Value = 22300
$ws.Range("D5").Value = 19400
$ws.Range("E5").Value = 24200

# --- Step 4: rename sheet
$ws.Name = "GWC"
